$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple text/percentage updates (safe as literal text) ---
$ws.Range("D2").Value = '39.937.57'
$ws.Range("E2").Value = '  -0.53%  '
$ws.Range("D3").Value = '2.217.95'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("E6").Value = '  -1.90%  '
$ws.Range("E7").Value = '  -0.81%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("E11").Value = '  +5.65%  '
$ws.Range("E12").Value = '  -1.50%  '
$ws.Range("E13").Value = '  +3.27%  '
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").Value = '2.562.19'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("E16").Value = '  -3.06%  '
$ws.Range("D17").Value = '2.257.19'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("D19").Value = '39.870.60'
$ws.Range("E19").Value = '  -0.57%  '
$ws.Range("D20").Value = '0.0₃0885'
$ws.Range("E20").Value = '  -0.52%  '
$ws.Range("E21").Value = '  -4.43%  '
$ws.Range("E22").Value = '  -2.20%  '
$ws.Range("E23").Value = '  -0.93%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("E26").Value = '  -0.99%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("E31").Value = '  +3.00%  '
$ws.Range("E32").Value = '  -3.88%  '
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("E35").Value = '  +4.51%  '
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("E37").Value = '  -1.40%  '
$ws.Range("E38").Value = '  -0.76%  '
$ws.Range("E39").Value = '  -1.53%  '
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("E41").Value = '  -6.41%  '
$ws.Range("D42").Value = '2.087.49'
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("E43").Value = '  -3.92%  '
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("E46").Value = '  -3.22%  '
$ws.Range("E47").Value = '  -9.04%  '
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("D49").Value = '2.435.49'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("E50").Value = '  -0.87%  '
$ws.Range("E51").Value = '  +1.89%  '

# --- Price updates that look like plain numbers: force text so Excel
#     does not coerce them (and drop trailing zeros / decimals) ---
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}
Set-TextValue $ws.Range("D5") '291.71'
Set-TextValue $ws.Range("D6") '86.76'
Set-TextValue $ws.Range("D10") '30.41'
Set-TextValue $ws.Range("D11") '50.38'
Set-TextValue $ws.Range("D14") '6.42'
Set-TextValue $ws.Range("D18") '0.731'
Set-TextValue $ws.Range("D21") '11.11'
Set-TextValue $ws.Range("D22") '5.74'
Set-TextValue $ws.Range("D23") '65.50'
Set-TextValue $ws.Range("D24") '236.84'
Set-TextValue $ws.Range("D26") '2.46'
Set-TextValue $ws.Range("D31") '157.51'
Set-TextValue $ws.Range("D32") '31.75'
Set-TextValue $ws.Range("D33") '0.999'
Set-TextValue $ws.Range("D36") '0.0713'
Set-TextValue $ws.Range("D39") '0.0987'
Set-TextValue $ws.Range("D41") '15.17'
Set-TextValue $ws.Range("D45") '17.90'
Set-TextValue $ws.Range("D48") '2.69'

# --- Rows 28/29: EthereumClassic and Toncoin swap places ---
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range("D28") '2.37'
$ws.Range("E28").Value = '  +7.64%  '

$ws.Range("B29").Value = 'EthereumClassic'
$ws.Range("C29").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D29") '23.28'
$ws.Range("E29").Value = '  +1.28%  '
